$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Rules")

# Change cell E8 text from "Good Morning" to "GIT UPDATE"
$ws.Range("E8").Value = "GIT UPDATE"

# Set the selected cell/range on the active sheet to E8
$ws.Range("E8").Select()
